# Updates the crypto price/volume table with refreshed values.
# Note: several "Price" (column D) values look like plain numbers (e.g. 521.48)
# but must stay as literal text, matching the source data which stores them as
# strings (e.g. "58.643.73" uses dots as thousands separators, not a valid
# number). A leading apostrophe forces Excel to keep such values as text
# instead of auto-converting them to floating point numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '58.672.64'
$ws.Range("E2").Value = '  -0.92%  '
$ws.Range("D3").Value = '2.630.83'
$ws.Range("E3").Value = '  -0.29%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Value = '''521.48'
$ws.Range("E5").Value = '  +1.00%  '
$ws.Range("D6").Value = '''144.80'
$ws.Range("E6").Value = '  -3.22%  '
$ws.Range("E7").Value = '  +0.28%  '
$ws.Range("D8").Value = '''0.576'
$ws.Range("E8").Value = '  -0.41%  '
$ws.Range("D9").Value = '2.638.97'
$ws.Range("E9").Value = '  -1.00%  '
$ws.Range("D10").Value = '''6.26'
$ws.Range("E10").Value = '  -3.89%  '
$ws.Range("E11").Value = '  -2.62%  '
$ws.Range("E12").Value = '  -2.35%  '
$ws.Range("E13").Value = '  -0.68%  '
$ws.Range("D14").Value = '3.088.20'
$ws.Range("E14").Value = '  -0.32%  '
$ws.Range("D15").Value = '58.662.67'
$ws.Range("E15").Value = '  -0.62%  '
$ws.Range("D16").Value = '''20.75'
$ws.Range("E16").Value = '  -3.45%  '
$ws.Range("E17").Value = '  -2.87%  '
$ws.Range("D18").Value = '2.636.74'
$ws.Range("E18").Value = '  -0.61%  '
$ws.Range("D19").Value = '''346.14'
$ws.Range("E19").Value = '  -0.59%  '
$ws.Range("D20").Value = '''4.43'
$ws.Range("E20").Value = '  -4.34%  '
$ws.Range("D21").Value = '''10.20'
$ws.Range("E21").Value = '  -3.86%  '
$ws.Range("E22").Value = '  -2.51%  '
$ws.Range("D23").Value = '''0.999'
$ws.Range("E23").Value = '  +0.02%  '
$ws.Range("D24").Value = '''61.71'
$ws.Range("E24").Value = '  +0.87%  '
$ws.Range("E25").Value = '  -3.20%  '
$ws.Range("E26").Value = '  +0.37%  '
$ws.Range("E27").Value = '  +0.40%  '
$ws.Range("D28").Value = '0.0₃0800'
$ws.Range("E28").Value = '  -4.46%  '
$ws.Range("D29").Value = '''7.01'
$ws.Range("E29").Value = '  -2.19%  '
$ws.Range("E30").Value = '  +0.17%  '
$ws.Range("E31").Value = '  -3.65%  '
$ws.Range("E32").Value = '  -0.18%  '
$ws.Range("E33").Value = '  -2.05%  '
$ws.Range("D34").Value = '''149.13'
$ws.Range("E34").Value = '  -0.05%  '
$ws.Range("D35").Value = '''0.974'
$ws.Range("E35").Value = '  -7.32%  '
$ws.Range("E36").Value = '  -3.18%  '
$ws.Range("E37").Value = '  -2.45%  '
$ws.Range("D38").Value = '''36.57'
$ws.Range("E38").Value = '  +0.54%  '
$ws.Range("D39").Value = '''0.836'
$ws.Range("E39").Value = '  -6.46%  '
$ws.Range("E40").Value = '  -2.82%  '
$ws.Range("D41").Value = '''3.61'
$ws.Range("E41").Value = '  -2.21%  '
$ws.Range("D42").Value = '''279.95'
$ws.Range("E42").Value = '  -5.29%  '
$ws.Range("E43").Value = '  +0.57%  '
$ws.Range("D44").Value = '''0.0984'
$ws.Range("E44").Value = '  -1.74%  '
$ws.Range("B45").Value = 'Mantle'
$ws.Range("C45").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D45").Value = '''0.598'
$ws.Range("E45").Value = '  -4.70%  '
$ws.Range("B46").Value = 'EnergySwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D46").Value = '''19.49'
$ws.Range("E46").Value = '  -2.16%  '
$ws.Range("E47").Value = '  -4.86%  '
$ws.Range("E48").Value = '  +0.39%  '
$ws.Range("E49").Value = '  -2.50%  '
$ws.Range("D50").Value = '1.984.00'
$ws.Range("E50").Value = '  +0.20%  '
$ws.Range("D51").Value = '''4.62'
$ws.Range("E51").Value = '  -3.63%  '
